$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 4 values per diff
$ws.Range("G4").Value = 1.83
$ws.Range("I4").Value = 5.5
$ws.Range("J4").Value = 2.63
$ws.Range("L4").Value = 6
$ws.Range("M4").Value = 1.14
$ws.Range("N4").Value = 5.5
$ws.Range("AC4").Value = 5.5
$ws.Range("AH4").Value = 9
$ws.Range("AI4").Value = 23
$ws.Range("AM4").Value = 67
$ws.Range("AW4").Value = 6.5

# Delete row 8 entirely (shift rows up)
$ws.Rows.Item(8).Delete()
